# Update gh-pages to output generated at 456a3b4
#
# Sheet layout:
#   1 = 展览 (Exhibition)
#   2 = 演出 (Performance)
#   3 = 本地生活 (Local Life)
#   4 = 全部类型 (All Types)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a literal text value into a cell without letting Excel's
# autodetect coerce date-looking strings ("2024-07-27") into date
# serials. We flip the cell to text format first, assign the value, then
# reset the style back to Normal so no stray formatting is left behind.
# ---------------------------------------------------------------------
function Set-TextValue {
    param($cell, [string]$value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# =======================================================================
# Sheet 1 - 展览: bump "想去人数" (interest count) in column F
# =======================================================================
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(3,6).Value  = 316
$ws1.Cells.Item(4,6).Value  = 2973
$ws1.Cells.Item(7,6).Value  = 2319
$ws1.Cells.Item(8,6).Value  = 1671
$ws1.Cells.Item(10,6).Value = 852
$ws1.Cells.Item(11,6).Value = 120
$ws1.Cells.Item(13,6).Value = 2662
$ws1.Cells.Item(15,6).Value = 1520
$ws1.Cells.Item(16,6).Value = 7056
$ws1.Cells.Item(18,6).Value = 7212
$ws1.Cells.Item(21,6).Value = 5452
$ws1.Cells.Item(23,6).Value = 3477
$ws1.Cells.Item(25,6).Value = 177
$ws1.Cells.Item(26,6).Value = 1882
$ws1.Cells.Item(33,6).Value = 2417
$ws1.Cells.Item(34,6).Value = 1184
$ws1.Cells.Item(35,6).Value = 2693
$ws1.Cells.Item(36,6).Value = 27
$ws1.Cells.Item(40,6).Value = 1074
$ws1.Cells.Item(42,6).Value = 478
$ws1.Cells.Item(43,6).Value = 524

# =======================================================================
# Sheet 2 - 演出: bump a few counts, then insert a brand-new row 19 for
# the "童年时光机" concert (pushing the old row 19 down to row 20).
# =======================================================================
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(12,6).Value = 43
$ws2.Cells.Item(16,6).Value = 18
$ws2.Cells.Item(18,6).Value = 62

# Insert a new row at position 19; existing row 19 (and below) shifts to 20.
$ws2.Rows.Item(19).Insert()

# column A keeps the bold/bordered header style used by every other row;
# Insert() alone doesn't carry it over, so copy it explicitly from A1.
$ws2.Range("A1").Copy() | Out-Null
$ws2.Cells.Item(19,1).PasteSpecial(-4122) | Out-Null

# --- new row 19: 北京·"童年时光机" concert ---
$ws2.Cells.Item(19,1).Value = 18
Set-TextValue $ws2.Cells.Item(19,2) "2024-07-27"
Set-TextValue $ws2.Cells.Item(19,3) "北京·“童年时光机”——《哆啦A梦》、《灌篮高手》、《狮子王》致敬童年经典动漫交响音乐会"
Set-TextValue $ws2.Cells.Item(19,4) "中关村南大街33号国家图书馆北门 国图艺术中心音乐厅"
Set-TextValue $ws2.Cells.Item(19,5) "2024.07.27 19:30-07.27 21:00"
$ws2.Cells.Item(19,6).Value = 0
Set-TextValue $ws2.Cells.Item(19,7) "不可售"
Set-TextValue $ws2.Cells.Item(19,8) "https://show.bilibili.com/platform/detail.html?id=85671"
Set-TextValue $ws2.Cells.Item(19,9) "//i2.hdslb.com/bfs/openplatform/202405/KV93ax2g1715669330587.jpeg"

# --- row 20 (formerly row 19, shifted down): fix the running index A20 ---
$ws2.Cells.Item(20,1).Value = 19

# =======================================================================
# Sheet 3 - 本地生活: add the first-ever data row (row 2), the sheet was
# previously just the header row.
# =======================================================================
$ws3 = $wb.Worksheets.Item(3)

# column A keeps the bold/bordered header style used on row 1.
$ws3.Range("A1").Copy() | Out-Null
$ws3.Cells.Item(2,1).PasteSpecial(-4122) | Out-Null

$ws3.Cells.Item(2,1).Value = 1
Set-TextValue $ws3.Cells.Item(2,2) "2024-05-25"
Set-TextValue $ws3.Cells.Item(2,3) "北京·战双帕弥什 x HAPPY ZOO主题Cafe"
Set-TextValue $ws3.Cells.Item(2,4) "学清路38号金码大厦B座(六道口地铁站B东北口步行110米) BOM嘻番里"
Set-TextValue $ws3.Cells.Item(2,5) "2024.05.25 00:00-06.03 23:59"
$ws3.Cells.Item(2,6).Value = 15
$ws3.Cells.Item(2,7).Value = 10
Set-TextValue $ws3.Cells.Item(2,8) "https://show.bilibili.com/platform/detail.html?id=85652"
Set-TextValue $ws3.Cells.Item(2,9) "//i0.hdslb.com/bfs/openplatform/202405/yVUhCFNH1715760749337.png"

# =======================================================================
# Sheet 4 - 全部类型: same "想去人数" bumps as sheet 1, mirrored into this
# combined view (this sheet keeps its original row layout - no inserts).
# =======================================================================
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(3,6).Value  = 316
$ws4.Cells.Item(5,6).Value  = 2973
$ws4.Cells.Item(7,6).Value  = 2319
$ws4.Cells.Item(8,6).Value  = 1671
$ws4.Cells.Item(11,6).Value = 852
$ws4.Cells.Item(12,6).Value = 120
$ws4.Cells.Item(14,6).Value = 2662
$ws4.Cells.Item(15,6).Value = 1520
$ws4.Cells.Item(19,6).Value = 7056
$ws4.Cells.Item(21,6).Value = 7212
$ws4.Cells.Item(23,6).Value = 5452
$ws4.Cells.Item(25,6).Value = 3477
$ws4.Cells.Item(29,6).Value = 1882
$ws4.Cells.Item(30,6).Value = 18
$ws4.Cells.Item(37,6).Value = 2417
$ws4.Cells.Item(38,6).Value = 1184
$ws4.Cells.Item(39,6).Value = 62
$ws4.Cells.Item(40,6).Value = 2693
$ws4.Cells.Item(41,6).Value = 27
$ws4.Cells.Item(46,6).Value = 1074
$ws4.Cells.Item(48,6).Value = 478
$ws4.Cells.Item(49,6).Value = 524
